$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header row): extend with P1=14, Q1=15, matching style of O1 ---
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows 2-25: update I/K/M/O values and add P/Q columns ---
for ($r = 2; $r -le 25; $r++) {
    # extend formatting of the new P/Q cells to match the rest of the (unstyled) data row
    $ws.Range("N$r").Copy() | Out-Null
    $ws.Range("P$r`:Q$r").PasteSpecial(-4122) | Out-Null

    $ws.Range("I$r").Value = 2
    $ws.Range("K$r").Value = 1
    $ws.Range("M$r").Value = 2
    $ws.Range("O$r").Value = 1
    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 2
}
